$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for Price column cells being updated, so numeric-looking
# strings (e.g. "298.09", "0.490") are preserved exactly as text, matching
# the original inlineStr cell type rather than being parsed as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = "42.125.90"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.265.72"
$ws.Range("E3").Value = "  -3.09%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "298.09"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "94.22"
$ws.Range("E6").Value = "  -7.01%  "
$ws.Range("E7").Value = "  -2.88%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  -4.15%  "
$ws.Range("D10").Value = "32.99"
$ws.Range("E10").Value = "  -5.55%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "48.19"
$ws.Range("E12").Value = "  -7.78%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "15.54"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.616.54"
$ws.Range("E16").Value = "  -3.21%  "
$ws.Range("D17").Value = "2.251.30"
$ws.Range("E17").Value = "  -4.11%  "
$ws.Range("D18").Value = "0.773"
$ws.Range("E18").Value = "  -5.40%  "
$ws.Range("D19").Value = "42.097.45"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").Value = "0.0₃0889"
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "5.98"
$ws.Range("E21").Value = "  -4.00%  "
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  -3.66%  "
$ws.Range("E23").Value = "  -1.88%  "
$ws.Range("D24").Value = "233.20"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").Value = "23.80"
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").Value = "167.64"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "33.68"
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("E32").Value = "  -4.20%  "
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "4.50"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D37").Value = "0.0690"
$ws.Range("E37").Value = "  -5.45%  "
$ws.Range("D38").Value = "16.21"
$ws.Range("E38").Value = "  -7.98%  "
$ws.Range("E39").Value = "  -5.17%  "
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  -8.30%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "1.960.47"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("E45").Value = "  -3.19%  "
$ws.Range("D46").Value = "17.33"
$ws.Range("E46").Value = "  -7.83%  "
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -5.52%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "2.80"
$ws.Range("E49").Value = "  -3.22%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.489.92"
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").Value = "51.85"
$ws.Range("E51").Value = "  -7.17%  "
